$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -eq "selfhelp") {
        $cell.Value = "unguided"
    } elseif ($val -eq "support") {
        $cell.Value = "guided"
    }
}
